$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 128.93333
$ws.Range("I33").Value = 131
$ws.Range("K33").Value = 131
$ws.Range("M33").Value = 98

$ws.Range("H55").Value = 232.4
$ws.Range("J55").Value = 289.85715
$ws.Range("L55").Value = 289.85715
$ws.Range("N55").Value = -717.85715

$ws.Range("H74").Value = 3599
$ws.Range("I74").Value = 3599
$ws.Range("K74").Value = 3599
$ws.Range("M74").Value = -2663

$ws.Range("H77").Value = 3599
$ws.Range("I77").Value = 3599
$ws.Range("K77").Value = 17995
$ws.Range("M77").Value = -13315

$ws.Range("H100").Value = 1381.4
$ws.Range("I100").Value = 851.875
$ws.Range("K100").Value = 851.875
$ws.Range("M100").Value = -310.875

$ws.Range("H106").Value = 34616.445
$ws.Range("I106").Value = 35818.5
$ws.Range("J106").Value = 25000
$ws.Range("K106").Value = 35818.5
$ws.Range("L106").Value = 25000
$ws.Range("M106").Value = -35187.5
$ws.Range("N106").Value = -26262

$ws.Range("H137").Value = 3082.577
$ws.Range("I137").Value = 1427.0834
$ws.Range("K137").Value = 4281.2502
$ws.Range("M137").Value = -1731.2502

$ws.Range("H141").Value = 5069.5835
$ws.Range("I141").Value = 4439.5454
$ws.Range("K141").Value = 13318.6362
$ws.Range("M141").Value = -8138.636200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7130.7554
$ws.Range("I32").Value = 5387.4146
$ws.Range("K32").Value = 5387.4146
$ws.Range("M32").Value = -5100.4146

$ws.Range("H55").Value = 35999.5
$ws.Range("J55").Value = 69999
$ws.Range("L55").Value = 69999
$ws.Range("N55").Value = -70629

$ws.Range("H61").Value = 4473
$ws.Range("I61").Value = 4122.3335
$ws.Range("K61").Value = 4122.3335
$ws.Range("M61").Value = -3910.3335

$ws.Range("H102").Value = 1372.5454
$ws.Range("I102").Value = 1199.75
$ws.Range("J102").Value = 1833.3334
$ws.Range("K102").Value = 1199.75
$ws.Range("L102").Value = 1833.3334
$ws.Range("M102").Value = 422.25
$ws.Range("N102").Value = -5077.3334

$ws.Range("H122").Value = 3493.6428
$ws.Range("I122").Value = 2582.625
$ws.Range("J122").Value = 4708.3335
$ws.Range("K122").Value = 7747.875
$ws.Range("L122").Value = 14125.0005
$ws.Range("M122").Value = -5297.875
$ws.Range("N122").Value = -19025.0005

$ws.Range("H132").Value = 1720.6945
$ws.Range("I132").Value = 1474.875
$ws.Range("J132").Value = 3687.25
$ws.Range("K132").Value = 4424.625
$ws.Range("L132").Value = 11061.75
$ws.Range("M132").Value = -1894.625
$ws.Range("N132").Value = -16121.75

$ws.Range("H136").Value = 4473
$ws.Range("I136").Value = 4122.3335
$ws.Range("K136").Value = 12367.0005
$ws.Range("M136").Value = -9817.000499999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 328.07144
$ws.Range("J80").Value = 176.6
$ws.Range("L80").Value = 176.6
$ws.Range("N80").Value = -2172.6

$ws.Range("H83").Value = 328.07144
$ws.Range("J83").Value = 176.6
$ws.Range("L83").Value = 883
$ws.Range("N83").Value = -10867

$ws.Range("H105").Value = 3351.9656
$ws.Range("I105").Value = 2455.8518
$ws.Range("K105").Value = 2455.8518
$ws.Range("M105").Value = -708.8517999999999

$ws.Range("H107").Value = 712.86365
$ws.Range("I107").Value = 714.9474
$ws.Range("K107").Value = 714.9474
$ws.Range("M107").Value = 1205.0526

$ws.Range("H134").Value = 2897.5833
$ws.Range("I134").Value = 2548.875
$ws.Range("J134").Value = 3595
$ws.Range("K134").Value = 7646.625
$ws.Range("L134").Value = 10785
$ws.Range("M134").Value = -5111.625
$ws.Range("N134").Value = -15855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 666.35297
$ws.Range("I16").Value = 642.75
$ws.Range("K16").Value = 642.75
$ws.Range("M16").Value = -355.75

$ws.Range("H31").Value = 4835.6
$ws.Range("I31").Value = 2819.6667
$ws.Range("K31").Value = 2819.6667
$ws.Range("M31").Value = -2524.6667

$ws.Range("H34").Value = 4835.6
$ws.Range("I34").Value = 2819.6667
$ws.Range("K34").Value = 2819.6667
$ws.Range("M34").Value = -2617.6667

$ws.Range("H52").Value = 87849.5
$ws.Range("J52").Value = 87849.5
$ws.Range("L52").Value = 87849.5
$ws.Range("N52").Value = -88437.5

$ws.Range("H113").Value = 666.35297
$ws.Range("I113").Value = 642.75
$ws.Range("K113").Value = 642.75
$ws.Range("M113").Value = 1527.25

$ws.Range("H134").Value = 3030.3215
$ws.Range("I134").Value = 2215.1667
$ws.Range("K134").Value = 6645.500100000001
$ws.Range("M134").Value = -4110.500100000001

$ws.Range("H138").Value = 99000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 473.90475
$ws.Range("J5").Value = 468.57144
$ws.Range("L5").Value = 1405.71432
$ws.Range("N5").Value = -1629.71432

$ws.Range("H68").Value = 1452.75
$ws.Range("I68").Value = 1206.2858
$ws.Range("K68").Value = 3618.8574
$ws.Range("M68").Value = -2807.8574

$ws.Range("H71").Value = 1452.75
$ws.Range("I71").Value = 1206.2858
$ws.Range("K71").Value = 10856.5722
$ws.Range("M71").Value = -6800.572200000001

$ws.Range("H135").Value = 473.90475
$ws.Range("J135").Value = 468.57144
$ws.Range("L135").Value = 4217.14296
$ws.Range("N135").Value = -9287.142960000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 651.2273
$ws.Range("I97").Value = 495.94116
$ws.Range("J97").Value = 1179.2
$ws.Range("K97").Value = 495.94116
$ws.Range("L97").Value = 1179.2
$ws.Range("M97").Value = 0.05883999999997513
$ws.Range("N97").Value = -2171.2

$ws.Range("H107").Value = 679.5454999999999
$ws.Range("I107").Value = 322.44446
$ws.Range("J107").Value = 926.7692
$ws.Range("K107").Value = 322.44446
$ws.Range("L107").Value = 926.7692
$ws.Range("M107").Value = 1597.55554
$ws.Range("N107").Value = -4766.7692

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1814.6666
$ws.Range("I40").Value = 1814.6666
$ws.Range("K40").Value = 1814.6666
$ws.Range("M40").Value = -1678.6666

$ws.Range("H46").Value = 1182.8334
$ws.Range("I46").Value = 1165.6666
$ws.Range("J46").Value = 1200
$ws.Range("K46").Value = 1165.6666
$ws.Range("L46").Value = 1200
$ws.Range("M46").Value = -977.6666
$ws.Range("N46").Value = -1576

$ws.Range("H93").Value = 2895.5454
$ws.Range("I93").Value = 2895.5454
$ws.Range("K93").Value = 2895.5454
$ws.Range("M93").Value = -1647.5454

$ws.Range("H132").Value = 4242.6665
$ws.Range("I132").Value = 4489.25
$ws.Range("J132").Value = 3749.5
$ws.Range("K132").Value = 13467.75
$ws.Range("L132").Value = 11248.5
$ws.Range("M132").Value = -10937.75
$ws.Range("N132").Value = -16308.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H111").Value = 80000
$ws.Range("J111").Value = 80000
$ws.Range("L111").Value = 80000
$ws.Range("N111").Value = -88180

$ws.Range("H126").Value = 2077.1052
$ws.Range("I126").Value = 1497.1875
$ws.Range("K126").Value = 4491.5625
$ws.Range("M126").Value = -2021.5625

$ws.Range("H136").Value = 5529.3335
$ws.Range("I136").Value = 5295
$ws.Range("J136").Value = 5998
$ws.Range("K136").Value = 15885
$ws.Range("L136").Value = 17994
$ws.Range("M136").Value = -13335
$ws.Range("N136").Value = -23094

$ws.Range("H141").Value = 50000
$ws.Range("I141").Value = 50000
$ws.Range("K141").Value = 50000
$ws.Range("M141").Value = -44820
